$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 54.69462833333333
$ws.Range("H2").Value = 164.083885
$ws.Range("I2").Value = 0.2790924419198448
$ws.Range("J2").Value = 0.2790924419198448
$ws.Range("M2").Value = 28.689524
$ws.Range("N2").Value = 86.06857199999999
$ws.Range("O2").Value = 0.2394085694101769
$ws.Range("P2").Value = 0.2394085694101769
$ws.Range("Q2").Value = 1569.162852240246
$ws.Range("R2").Value = 14122.46567016222
$ws.Range("S2").Value = 0.06681712225322294
$ws.Range("T2").Value = 0.06681712225322294
$ws.Range("G3").Value = 54.69462833333333
$ws.Range("H3").Value = 164.083885
$ws.Range("I3").Value = 0.2790924419198448
$ws.Range("J3").Value = 0.2790924419198448
$ws.Range("O3").Value = 0.5212694246546397
$ws.Range("P3").Value = 0.5212694246546395
$ws.Range("Q3").Value = 3416.572009898724
$ws.Range("R3").Value = 30749.14808908852
$ws.Range("S3").Value = 0.1454823566250159
$ws.Range("T3").Value = 0.1454823566250159
$ws.Range("G4").Value = 54.69462833333333
$ws.Range("H4").Value = 164.083885
$ws.Range("I4").Value = 0.2790924419198448
$ws.Range("J4").Value = 0.2790924419198448
$ws.Range("M4").Value = 28.525746
$ws.Range("N4").Value = 85.577238
$ws.Range("O4").Value = 0.2380418734454457
$ws.Range("P4").Value = 0.2380418734454457
$ws.Range("Q4").Value = 1560.20507540107
$ws.Range("R4").Value = 14041.84567860963
$ws.Range("S4").Value = 0.06643568773906411
$ws.Range("T4").Value = 0.0664356877390641
$ws.Range("G5").Value = 54.69462833333333
$ws.Range("H5").Value = 164.083885
$ws.Range("I5").Value = 0.2790924419198448
$ws.Range("J5").Value = 0.2790924419198448
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 0.1534046666666667
$ws.Range("N5").Value = 0.460214
$ws.Range("O5").Value = 0.001280132489737778
$ws.Range("P5").Value = 0.001280132489737778
$ws.Range("Q5").Value = 8.390411227932223
$ws.Range("R5").Value = 75.51370105139
$ws.Range("S5").Value = 0.000357275302541847
$ws.Range("T5").Value = 0.000357275302541847
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09859081282432611
$ws.Range("J6").Value = 0.09859081282432611
$ws.Range("M6").Value = 28.689524
$ws.Range("N6").Value = 86.06857199999999
$ws.Range("O6").Value = 0.2394085694101769
$ws.Range("P6").Value = 0.2394085694101769
$ws.Range("Q6").Value = 554.3146922643467
$ws.Range("R6").Value = 4988.83223037912
$ws.Range("S6").Value = 0.02360348545525844
$ws.Range("T6").Value = 0.02360348545525844
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09859081282432611
$ws.Range("J7").Value = 0.09859081282432611
$ws.Range("O7").Value = 0.5212694246546397
$ws.Range("P7").Value = 0.5212694246546395
$ws.Range("S7").Value = 0.05139237627716974
$ws.Range("T7").Value = 0.05139237627716973
$ws.Range("G8").Value = 19.32115333333334
$ws.Range("H8").Value = 57.96346000000001
$ws.Range("I8").Value = 0.09859081282432611
$ws.Range("J8").Value = 0.09859081282432611
$ws.Range("M8").Value = 28.525746
$ws.Range("N8").Value = 85.577238
$ws.Range("O8").Value = 0.2380418734454457
$ws.Range("P8").Value = 0.2380418734454457
$ws.Range("Q8").Value = 551.1503124137201
$ws.Range("R8").Value = 4960.352811723481
$ws.Range("S8").Value = 0.02346874178921186
$ws.Range("T8").Value = 0.02346874178921186
$ws.Range("G9").Value = 19.32115333333334
$ws.Range("H9").Value = 57.96346000000001
$ws.Range("I9").Value = 0.09859081282432611
$ws.Range("J9").Value = 0.09859081282432611
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 0.1534046666666667
$ws.Range("N9").Value = 0.460214
$ws.Range("O9").Value = 0.001280132489737778
$ws.Range("P9").Value = 0.001280132489737778
$ws.Range("Q9").Value = 2.963955086715556
$ws.Range("R9").Value = 26.67559578044001
$ws.Range("S9").Value = 0.0001262093026860758
$ws.Range("T9").Value = 0.0001262093026860758
$ws.Range("G10").Value = 11.023718
$ws.Range("H10").Value = 33.071154
$ws.Range("I10").Value = 0.05625116157486912
$ws.Range("J10").Value = 0.05625116157486911
$ws.Range("M10").Value = 28.689524
$ws.Range("N10").Value = 86.06857199999999
$ws.Range("O10").Value = 0.2394085694101769
$ws.Range("P10").Value = 0.2394085694101769
$ws.Range("Q10").Value = 316.265222130232
$ws.Range("R10").Value = 2846.386999172088
$ws.Range("S10").Value = 0.01346701012030013
$ws.Range("T10").Value = 0.01346701012030013
$ws.Range("G11").Value = 11.023718
$ws.Range("H11").Value = 33.071154
$ws.Range("I11").Value = 0.05625116157486912
$ws.Range("J11").Value = 0.05625116157486911
$ws.Range("O11").Value = 0.5212694246546397
$ws.Range("P11").Value = 0.5212694246546395
$ws.Range("Q11").Value = 688.6110667811787
$ws.Range("R11").Value = 6197.499601030609
$ws.Range("S11").Value = 0.0293220106302872
$ws.Range("T11").Value = 0.02932201063028719
$ws.Range("G12").Value = 11.023718
$ws.Range("H12").Value = 33.071154
$ws.Range("I12").Value = 0.05625116157486912
$ws.Range("J12").Value = 0.05625116157486911
$ws.Range("M12").Value = 28.525746
$ws.Range("N12").Value = 85.577238
$ws.Range("O12").Value = 0.2380418734454457
$ws.Range("P12").Value = 0.2380418734454457
$ws.Range("Q12").Value = 314.459779643628
$ws.Range("R12").Value = 2830.138016792652
$ws.Range("S12").Value = 0.01339013188476431
$ws.Range("T12").Value = 0.01339013188476431
$ws.Range("G13").Value = 11.023718
$ws.Range("H13").Value = 33.071154
$ws.Range("I13").Value = 0.05625116157486912
$ws.Range("J13").Value = 0.05625116157486911
$ws.Range("K13").Value = 3.0
$ws.Range("L13").Value = 1.0
$ws.Range("M13").Value = 0.1534046666666667
$ws.Range("N13").Value = 0.460214
$ws.Range("O13").Value = 0.001280132489737778
$ws.Range("P13").Value = 0.001280132489737778
$ws.Range("Q13").Value = 1.691089785217333
$ws.Range("R13").Value = 15.219808066956
$ws.Range("S13").Value = 0.00007200893951747922
$ws.Range("T13").Value = 0.0000720089395174792
$ws.Range("G14").Value = 110.9336623333333
$ws.Range("H14").Value = 332.800987
$ws.Range("I14").Value = 0.5660655836809599
$ws.Range("J14").Value = 0.5660655836809599
$ws.Range("M14").Value = 28.689524
$ws.Range("N14").Value = 86.06857199999999
$ws.Range("O14").Value = 0.2394085694101769
$ws.Range("P14").Value = 0.2394085694101769
$ws.Range("Q14").Value = 3182.633967920061
$ws.Range("R14").Value = 28643.70571128056
$ws.Range("S14").Value = 0.1355209515813954
$ws.Range("T14").Value = 0.1355209515813954
$ws.Range("G15").Value = 110.9336623333333
$ws.Range("H15").Value = 332.800987
$ws.Range("I15").Value = 0.5660655836809599
$ws.Range("J15").Value = 0.5660655836809599
$ws.Range("O15").Value = 0.5212694246546397
$ws.Range("P15").Value = 0.5212694246546395
$ws.Range("Q15").Value = 6929.617354262847
$ws.Range("R15").Value = 62366.55618836562
$ws.Range("S15").Value = 0.2950726811221668
$ws.Range("T15").Value = 0.2950726811221667
$ws.Range("G16").Value = 110.9336623333333
$ws.Range("H16").Value = 332.800987
$ws.Range("I16").Value = 0.5660655836809599
$ws.Range("J16").Value = 0.5660655836809599
$ws.Range("M16").Value = 28.525746
$ws.Range("N16").Value = 85.577238
$ws.Range("O16").Value = 0.2380418734454457
$ws.Range("P16").Value = 0.2380418734454457
$ws.Range("Q16").Value = 3164.465474570433
$ws.Range("R16").Value = 28480.1892711339
$ws.Range("S16").Value = 0.1347473120324054
$ws.Range("T16").Value = 0.1347473120324054
$ws.Range("G17").Value = 110.9336623333333
$ws.Range("H17").Value = 332.800987
$ws.Range("I17").Value = 0.5660655836809599
$ws.Range("J17").Value = 0.5660655836809599
$ws.Range("K17").Value = 3.0
$ws.Range("L17").Value = 1.0
$ws.Range("M17").Value = 0.1534046666666667
$ws.Range("N17").Value = 0.460214
$ws.Range("O17").Value = 0.001280132489737778
$ws.Range("P17").Value = 0.001280132489737778
$ws.Range("Q17").Value = 17.01774149235755
$ws.Range("R17").Value = 153.159673431218
$ws.Range("S17").Value = 0.0007246389449923756
$ws.Range("T17").Value = 0.0007246389449923756
